# Update the text of the "NLP Prep Process" SmartArt diagram on slide 1.
# The diagram (a Step-Down Process SmartArt) is the single shape on the
# slide; its step titles/descriptions are rewritten per the authored edit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(1)
$sa = $shp.SmartArt
$nodes = $sa.AllNodes

# Step 1: LOAD
$nodes.Item(1).TextFrame2.TextRange.Text = "Load"
# (node 2 "Load corpus into an RDD" is unchanged)

# Step 2: LOWER -> Create / RDD (two-line title)
$nodes.Item(3).TextFrame2.TextRange.Text = "Create" + [char]13 + "RDD"
$nodes.Item(4).TextFrame2.TextRange.Text = "Define and create the RDD"

# Step 3: FILTER -> Flatten
$nodes.Item(5).TextFrame2.TextRange.Text = "Flatten"
$nodes.Item(6).TextFrame2.TextRange.Text = "Transform data frame into a flat structure"

# Step 4: SENTENCE TOKENIZE -> Remove Header
$nodes.Item(7).TextFrame2.TextRange.Text = "Remove Header"
$nodes.Item(8).TextFrame2.TextRange.Text = "Get rid of the header row"

# Step 5: WORDS TOKENIZE -> Sentence Tokenize
$nodes.Item(9).TextFrame2.TextRange.Text = "Sentence Tokenize"
$nodes.Item(10).TextFrame2.TextRange.Text = "Tokenize each sentence"

# Step 6: PUNCTUATION -> HTML Parsing
$nodes.Item(11).TextFrame2.TextRange.Text = "HTML Parsing"
$nodes.Item(12).TextFrame2.TextRange.Text = "Remove all HTML formats (intended for comments)"

# Step 7: LAMMETIZATION -> VADER
$nodes.Item(13).TextFrame2.TextRange.Text = "VADER"
$nodes.Item(14).TextFrame2.TextRange.Text = "Apply VADER sentiment analysis to each comment/headline"

# Step 8: RE-JOIN -> Transform to Pandas Data Frame
$nodes.Item(15).TextFrame2.TextRange.Text = "Transform to Pandas Data Frame"
$nodes.Item(16).TextFrame2.TextRange.Text = "Convert Spark data frame to Pandas data frame"
